# Insert a new data row at row 541 (pushing existing rows 541-591 down to
# 542-592) on the "Hortaliza, Femacal de La Calera - Cilantro" sheet, and
# populate the newly inserted row with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 541:591 down by one row to make room for the new observation.
$ws.Rows.Item(541).Insert()

# Fill in the new row 541 with the new record's data.
$ws.Cells.Item(541, 1).Value  = 3
$ws.Cells.Item(541, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(541, 3).Value  = "Coquimbo"
$ws.Cells.Item(541, 4).Value  = 45106
$ws.Cells.Item(541, 5).Value  = 5
$ws.Cells.Item(541, 6).Value  = 100112040
$ws.Cells.Item(541, 7).Value  = "Cilantro"
$ws.Cells.Item(541, 8).Value  = "Sin especificar"
$ws.Cells.Item(541, 9).Value  = "Primera"
$ws.Cells.Item(541, 10).Value = 240
$ws.Cells.Item(541, 11).Value = 3500
$ws.Cells.Item(541, 12).Value = 4000
$ws.Cells.Item(541, 13).Value = 3750
$ws.Cells.Item(541, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(541, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(541, 16).Value = 1250
$ws.Cells.Item(541, 17).Value = 3
$ws.Cells.Item(541, 18).Value = "Hortaliza"
